$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.344.04"
$ws.Range("E2").Value = "  -0.58%  "

$ws.Range("D3").Value = "'1.811.47"
$ws.Range("E3").Value = "  -0.88%  "

$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  -0.51%  "

$ws.Range("D5").Value = "'313.37"
$ws.Range("E5").Value = "  -0.98%  "

$ws.Range("D6").Value = "'0.9995"
$ws.Range("E6").Value = "  -0.44%  "

$ws.Range("D7").Value = "'0.5155"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'0.4008"
$ws.Range("E8").Value = "  +3.55%  "

$ws.Range("D9").Value = "'0.07874"
$ws.Range("E9").Value = "  -5.28%  "

$ws.Range("D10").Value = "'1.117"
$ws.Range("E10").Value = "  -0.27%  "

$ws.Range("D11").Value = "'41.06"
$ws.Range("E11").Value = "  -2.33%  "

$ws.Range("D12").Value = "'6.381"
$ws.Range("E12").Value = "  -0.57%  "

$ws.Range("D13").Value = "'0.9992"
$ws.Range("E13").Value = "  -0.51%  "

$ws.Range("D14").Value = "'20.48"
$ws.Range("E14").Value = "  -3.25%  "

$ws.Range("D15").Value = "'7.364"
$ws.Range("E15").Value = "  -1.84%  "

$ws.Range("D16").Value = "'1.799.61"
$ws.Range("E16").Value = "  -1.40%  "

$ws.Range("E17").Value = "  -1.20%  "

$ws.Range("D18").Value = "'0.00001084"
$ws.Range("E18").Value = "  -3.44%  "

$ws.Range("D19").Value = "'0.06588"
$ws.Range("E19").Value = "  -0.98%  "

$ws.Range("D20").Value = "'0.9989"
$ws.Range("E20").Value = "  -0.48%  "

$ws.Range("E21").Value = "  -2.51%  "

$ws.Range("D22").Value = "'6.042"
$ws.Range("E22").Value = "  -0.24%  "

$ws.Range("D23").Value = "'28.397.57"
$ws.Range("E23").Value = "  -0.54%  "

$ws.Range("D24").Value = "'11.22"
$ws.Range("E24").Value = "  -1.95%  "

$ws.Range("D25").Value = "'2.229"
$ws.Range("E25").Value = "  -2.56%  "

$ws.Range("D26").Value = "'160.73"
$ws.Range("E26").Value = "  +0.69%  "

$ws.Range("D27").Value = "'20.64"
$ws.Range("E27").Value = "  -2.52%  "

$ws.Range("D28").Value = "'2.016.80"
$ws.Range("E28").Value = "  -0.92%  "

$ws.Range("D29").Value = "'2.418"

$ws.Range("D30").Value = "'128.80"
$ws.Range("E30").Value = "  +2.09%  "

$ws.Range("D31").Value = "'0.1085"
$ws.Range("E31").Value = "  -0.72%  "

$ws.Range("D32").Value = "'1.054"
$ws.Range("E32").Value = "  -3.87%  "

$ws.Range("D33").Value = "'5.604"
$ws.Range("E33").Value = "  -2.39%  "

$ws.Range("D34").Value = "'3.662"
$ws.Range("E34").Value = "  -0.46%  "

$ws.Range("D35").Value = "'0.07203"
$ws.Range("E35").Value = "  -5.13%  "

$ws.Range("D36").Value = "'9.145"
$ws.Range("E36").Value = "  +4.36%  "

$ws.Range("D37").Value = "'0.02341"
$ws.Range("E37").Value = "  -1.67%  "

$ws.Range("D38").Value = "'0.2169"
$ws.Range("E38").Value = "  -2.82%  "

$ws.Range("E39").Value = "  -2.08%  "

$ws.Range("D40").Value = "'5.078"
$ws.Range("E40").Value = "  -3.77%  "

$ws.Range("D41").Value = "'0.6233"
$ws.Range("E41").Value = "  -2.25%  "

$ws.Range("D42").Value = "'0.9987"
$ws.Range("E42").Value = "  -0.49%  "

$ws.Range("D43").Value = "'1.159"
$ws.Range("E43").Value = "  -2.83%  "

$ws.Range("D44").Value = "'13.23"
$ws.Range("E44").Value = "  -2.77%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.6031"
$ws.Range("E45").Value = "  -1.24%  "

$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.313"
$ws.Range("E46").Value = "  -5.96%  "

$ws.Range("D47").Value = "'3.745"
$ws.Range("E47").Value = "  -1.58%  "

$ws.Range("E48").Value = "  -1.07%  "

$ws.Range("D49").Value = "'1.221"
$ws.Range("E49").Value = "  +1.19%  "

$ws.Range("D50").Value = "'1.946"
$ws.Range("E50").Value = "  -2.71%  "

$ws.Range("D51").Value = "'0.06856"
$ws.Range("E51").Value = "  -1.87%  "
